# The KeywordTok/ImportTok/.../ErrorTok character styles in styles.xml have
# their <w:rPr> children in the wrong order: wml.xsd (CT_RPr via the
# EG_RPrBase group) expects toggle properties such as <w:b/>/<w:i/> to
# precede <w:color/>, but these styles were serialized as
# <w:color/> followed by <w:b/>/<w:i/>. That trips OOXMLValidatorCLI's
# schema check (Sch_UnexpectedElementContentExpectingComplex) even though
# xmllint stays quiet.
#
# Re-assigning Font.Bold / Font.Italic to their own (already-true) value
# doesn't change any formatting, but it makes Word rewrite that style's
# <w:rPr> in schema order, moving <w:b/>/<w:i/> ahead of <w:color/>.
$d = $word.ActiveDocument

$boldOnly   = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
$italicOnly = @("CommentTok", "DocumentationTok")
$boldItalic = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")

foreach ($name in $boldOnly) {
    $d.Styles($name).Font.Bold = $true
}

foreach ($name in $italicOnly) {
    $d.Styles($name).Font.Italic = $true
}

foreach ($name in $boldItalic) {
    $d.Styles($name).Font.Bold = $true
    $d.Styles($name).Font.Italic = $true
}
